$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the values in A2 and A3 (finishing LP / OM used figures)
$ws.Range("A2").Value = 685601381667
$ws.Range("A3").Value = 685601355769

# Select A2:B3 with A2 as the active cell, matching the saved selection state
$ws.Range("A2:B3").Select()
$excel.ActiveCell = $ws.Range("A2")
